$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 97: A97 = 3, B97 = "নীয়"
$ws.Cells.Item(97, 1).Value = 3
$ws.Cells.Item(97, 2).Value = "নীয়"

# Copy formatting from row 96 to row 97 (so styles match the existing pattern)
$ws.Range("A96:B96").Copy()
$ws.Range("A97:B97").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to new last cell
$ws.Range("A97").Select()

# Adjust the window size (workbookView) to match the target
$excel.ActiveWindow.Width = 10920
$excel.ActiveWindow.Height = 9072
